# Change 1: split " file ALL of the animals/pics..." run into three runs
$d = $word.ActiveDocument

$r1 = $d.Content
$r1.Find.Execute("One idea, and this would be tedious", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r1.Expand(4) | Out-Null
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E219DA" w:rsidRDefault="00E219DA" w:rsidP="00E219DA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">One idea, and this would be tedious, would be to add to our “regular” resource file or some </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file </w:t></w:r><w:r><w:t xml:space="preserve">for </w:t></w:r><w:r><w:t xml:space="preserve">ALL of the animals/pics and what their translation is for each of the 9 possible languages.  For example, a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> setting for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LionFrench</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LionItallian</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LionChinese</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>...</w:t></w:r><w:r w:rsidR="001D33B4"><w:t xml:space="preserve">So if we have say 200 different pics, with 9 different languages, that would be 200*9, or 1800 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="001D33B4"><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="001D33B4"><w:t xml:space="preserve"> settings.</w:t></w:r></w:p>'
$r1.InsertXML($xml1)

# Change 2: insert new "TJY UPDATE" paragraph before "Need to introduce Pivot Pages",
# re-ordering the bookmark to the new paragraph.
$rStart = $d.Content
$rStart.Find.Execute("Need to introduce Pivot Pages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rStart.Expand(4) | Out-Null

$rEnd = $d.Content
$rEnd.Find.Execute("but they do not have enough)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rEnd.Expand(4) | Out-Null

$r2 = $d.Range($rStart.Start, $rEnd.End)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr><w:r><w:t>TJY UPDATE: For now the play/pause button is on the top next to the Baby Sounds app name.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EC0B37" w:rsidRDefault="00EC0B37" w:rsidP="00EC0B37"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Need to introduce Pivot Pages so user can move from page to page.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EC0B37" w:rsidRDefault="00EC0B37" w:rsidP="00EC0B37"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>App needs to be smart enough to not let user move to page if they only have trial offer.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008B4F62" w:rsidRDefault="00860515" w:rsidP="008B4F62"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>W</w:t></w:r><w:r w:rsidR="008B4F62"><w:t>e need more sounds and pics (below web sites are what I used, but they do not have enough)</w:t></w:r></w:p>'
$r2.InsertXML($xml2)

# Change 3: merge the two Consolas runs about pause/stop into a single run
$r3 = $d.Content
$r3.Find.Execute("FYI: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Expand(4) | Out-Null
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="001D33B4" w:rsidRPr="001D33B4" w:rsidRDefault="001D33B4" w:rsidP="001D33B4"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve">FYI: </w:t></w:r><w:r w:rsidRPr="001D33B4"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve">media element allows for pausing sound, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="001D33B4"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:highlight w:val="white"/></w:rPr><w:t>Soundeffect</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="001D33B4"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve"> does not allow for pause/stop so NOT good for background </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="001D33B4"><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:highlight w:val="white"/></w:rPr><w:t>musisc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$r3.InsertXML($xml3)

# Change 4: strike-through "Size of pics on grid needs to be uniform."
$r4 = $d.Content
$r4.Find.Execute("Size of pics on grid needs to be uniform.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Expand(4) | Out-Null
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EC0B37" w:rsidRDefault="00EC0B37" w:rsidP="00860515"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t>Size of pics on grid needs to be uniform.</w:t></w:r></w:p>'
$r4.InsertXML($xml4)

Write-Output "All changes applied"
